$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update SVM output note: it was still running, clarify as 1200+ minutes
# instead of a hard "1300 minutes".
$ws.Range("E13").Value = "Incomplete after 1200+ minutes"

# Fix the Random Forest parameter string for row 12 (it had an accidental
# duplicate "max_depth=150," before "max_depth=50"); also remove the now
# orphaned duplicate string from the shared string table implicitly by
# simply rewriting the cell value.
$ws.Range("E12").Value = "class_weight='balanced', criterion='entropy', max_depth=50, max_leaf_nodes=16, min_samples_leaf=6, min_samples_split=4, n_estimators=200"
